# "Analyse des performances apres 5eme correction"
# Fill in the ÉTAPE 5 (column H) results on all three report sheets,
# copying the formatting already used by ÉTAPE 4 (column G) on each row
# and writing in the step-5 score.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "LightHouse - Portable"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("G4").Copy()
$ws1.Range("H4").PasteSpecial(-4122)
$ws1.Range("H4").Value = 82

$ws1.Range("G5").Copy()
$ws1.Range("H5").PasteSpecial(-4122)
$ws1.Range("H5").Value = 88

$ws1.Range("G6").Copy()
$ws1.Range("H6").PasteSpecial(-4122)
$ws1.Range("H6").Value = 87

$ws1.Range("G7").Copy()
$ws1.Range("H7").PasteSpecial(-4122)
$ws1.Range("H7").Value = 78

$ws1.Range("H4:H7").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet 2: "LightHouse - Bureau"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("G4").Copy()
$ws2.Range("H4").PasteSpecial(-4122)
$ws2.Range("H4").Value = 92

$ws2.Range("G5").Copy()
$ws2.Range("H5").PasteSpecial(-4122)
$ws2.Range("H5").Value = 86

$ws2.Range("G6").Copy()
$ws2.Range("H6").PasteSpecial(-4122)
$ws2.Range("H6").Value = 93

$ws2.Range("G7").Copy()
$ws2.Range("H7").PasteSpecial(-4122)
$ws2.Range("H7").Value = 90

# Original selection was the non-contiguous "G4 G6 G7" -> "H4 H6:H7";
# this host's selection model only tracks a single rectangle, so the
# closest reachable state is the rectangle spanning the touched cells.
$ws2.Range("H4:H7").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet 3: "GTmetrix - Bureau"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("G4").Copy()
$ws3.Range("H4").PasteSpecial(-4122)
$ws3.Range("H4").Value = 100

$ws3.Range("G5").Copy()
$ws3.Range("H5").PasteSpecial(-4122)
$ws3.Range("H5").Value = 96

$ws3.Range("H4:H5").Select() | Out-Null
